$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2025-09-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-20 Saturday", 2) | Out-Null

# Update each cell of the single 20x5 table with the new equation text
$t = $d.Tables.Item(1)
$values = @(
    "31-20=11", "52+41=93", "3+4=7", "90-65=25", "81-4=77", "11+85=96", "30+0=30", "73+18=91",
    "91-88=3", "66+26=92", "3+90=93", "1+14=15", "21-8=13", "98-47=51", "99-16=83", "75-16=59",
    "33+62=95", "59+34=93", "45-31=14", "27+55=82", "6+16=22", "75-53=22", "54-43=11", "14+40=54",
    "59+18=77", "68-11=57", "86-25=61", "83-8=75", "66-32=34", "37+20=57", "14+81=95", "32+29=61",
    "18+52=70", "50+22=72", "0+55=55", "40-26=14", "94-28=66", "46+22=68", "87+0=87", "60+0=60",
    "20-4=16", "92-31=61", "27+4=31", "96+3=99", "19-18=1", "67-35=32", "1+71=72", "88+11=99",
    "99-77=22", "70-17=53", "22+76=98", "38+50=88", "85-76=9", "11+42=53", "84-20=64", "33+48=81",
    "58+36=94", "4+23=27", "93-27=66", "22+20=42", "51+32=83", "64-17=47", "54+45=99", "31-19=12",
    "89+10=99", "57-31=26", "69-22=47", "41-6=35", "59+36=95", "59-54=5", "99-80=19", "53-32=21",
    "12+55=67", "25+51=76", "5+48=53", "37+54=91", "50+1=51", "53+34=87", "21+7=28", "14+82=96",
    "21+48=69", "96+0=96", "70-18=52", "24+30=54", "17+81=98", "97-6=91", "8+5=13", "0+44=44",
    "96+2=98", "88-68=20", "9+77=86", "5+94=99", "16+45=61", "59-55=4", "66-20=46", "90-4=86",
    "67-60=7", "51+24=75", "53-37=16", "30-3=27"
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "done: $idx cells updated"
